$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 86, shifting existing rows 86-152 down to 87-153.
$ws.Rows.Item(86).Insert()

# Populate the newly inserted row 86 with the new weekly data point.
$ws.Range("A86").Value = 9
$ws.Range("B86").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C86").Value = "Metropolitana"
$ws.Range("D86").Value = 45271
$ws.Range("E86").Value = 13
$ws.Range("F86").Value = 100114007
$ws.Range("G86").Value = "Jengibre"
$ws.Range("H86").Value = "Sin especificar"
$ws.Range("I86").Value = "Primera"
$ws.Range("J86").Value = 340
$ws.Range("K86").Value = 19000
$ws.Range("L86").Value = 20000
$ws.Range("M86").Value = 19500
$ws.Range("N86").Value = "$/caja 13 kilos"
$ws.Range("O86").Value = "Perú"
$ws.Range("P86").Value = 1500
$ws.Range("Q86").Value = 13
$ws.Range("R86").Value = "Hortaliza"
